$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly observation at row 319 (pushes the existing
# 319..369 block down to 320..370, and R369 -> R370 in the used range).
$ws.Rows.Item(319).Insert()

$ws.Cells.Item(319, 1).Value = 7
$ws.Cells.Item(319, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(319, 3).Value = "Ñuble"
$ws.Cells.Item(319, 4).Value = 44951
$ws.Cells.Item(319, 5).Value = 16
$ws.Cells.Item(319, 6).Value = 100114013
$ws.Cells.Item(319, 7).Value = "Zanahoria"
$ws.Cells.Item(319, 8).Value = "Sin especificar"
$ws.Cells.Item(319, 9).Value = "Primera"
$ws.Cells.Item(319, 10).Value = 120
$ws.Cells.Item(319, 11).Value = 8500
$ws.Cells.Item(319, 12).Value = 9000
$ws.Cells.Item(319, 13).Value = 8750
$ws.Cells.Item(319, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(319, 15).Value = "Región de Ñuble"
$ws.Cells.Item(319, 16).Value = 438
$ws.Cells.Item(319, 17).Value = 20
$ws.Cells.Item(319, 18).Value = "Hortaliza"
